# Added soplex simulation results
# Fills in the previously-missing "soplex" benchmark rows (SRRIP / Hawkeye /
# OPTGen policies) on the Config1 sheet, and the "soplex" SRRIP row on the
# Config2 sheet, which were left blank (causing #DIV/0! in the dependent
# H/I formulas). Also updates the active sheet / selection to reflect where
# the author was working when the data was added.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Config1 ("soplex" rows 68-70): SRRIP, Hawkeye, OPTGen
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Config1")

# Row 68 - SRRIP
$ws1.Range("C68").Value = 50000000
$ws1.Range("D68").Value = 162827847
$ws1.Range("E68").Value = 2732027
$ws1.Range("F68").Value = 1497017
$ws1.Range("G68").Value = 1235010

# Row 69 - Hawkeye
$ws1.Range("C69").Value = 50000000
$ws1.Range("D69").Value = 156923926
$ws1.Range("E69").Value = 2731976
$ws1.Range("F69").Value = 1490724
$ws1.Range("G69").Value = 1241252

# Row 70 - OPTGen (G is derived E-F, plus a new J = F/E column)
$ws1.Range("C70").Value = 50000000
$ws1.Range("D70").Value = 156923926
$ws1.Range("E70").Value = 59591
$ws1.Range("F70").Value = 31306
$ws1.Range("G70").Formula = "=E70-F70"
$ws1.Range("J70").Formula = "=F70/E70"

# ----------------------------------------------------------------------
# Config2 ("soplex" row 68): SRRIP
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Config2")

$ws2.Range("C68").Value = 50000000
$ws2.Range("D68").Value = 124990062
$ws2.Range("E68").Value = 3394940
$ws2.Range("F68").Value = 1837023
$ws2.Range("G68").Value = 1557917

# ----------------------------------------------------------------------
# View state: the author ended up with Config1 active, scrolled to the
# newly-entered rows, and Config2 no longer the selected tab (its
# selection moved from C68 to C69, the row below the new data).
# ----------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("C69").Select() | Out-Null

$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 50
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("A50").Select() | Out-Null
